$wb = $excel.ActiveWorkbook

# "Generate Report for Handback": the shared "Ready for handoff" status text
# becomes "Handback transform failed", and the per-language sheets now carry
# a concrete error message in the "Error Detail" column.

$errorDetail = "The translationStateItem c8239dac9781b0284b70928afd51cfb0125338af is not found."
$failedStatus = "Handback transform failed"

# --- Overview sheet: zh-cn / de-de status cells ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $failedStatus
$wsOverview.Range("F2").Value = $failedStatus
# Widen the zh-cn / de-de columns to fit the longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 23.8333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 23.8333333333333

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $failedStatus
$wsZhCn.Range("P2").Value = $errorDetail
$wsZhCn.Columns.Item(3).ColumnWidth = 23.8333333333333
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1666666666667

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $failedStatus
$wsDeDe.Range("P2").Value = $errorDetail
$wsDeDe.Columns.Item(3).ColumnWidth = 23.8333333333333
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1666666666667
